$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New descriptions for RQ01..RQ10 (column B), replacing the old 14-row list
$descriptions = @(
    "Gestión de Autenticación y Autorización",
    "Registro de Usuarios",
    "Gestión de Establecimientos",
    "Sistema de Reservas",
    "Sistema de Convocatorias",
    "Búsqueda de Usuarios",
    "Administración de Usuarios",
    "Visualización de Reservas Administrativas",
    "Gestión de Perfiles",
    "Notificaciones por Email"
)

for ($i = 0; $i -lt $descriptions.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $descriptions[$i]
}

# Remove the now-obsolete rows 12-15 (RQ11-RQ14 and their descriptions)
$ws.Range("A12:B15").EntireRow.Delete() | Out-Null

# Update selection to match the resulting state after the edit
$ws.Range("A12:XFD15").Select() | Out-Null
